$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("B4")
$ws.Hyperlinks.Add($r, "https://example.com/test/") | Out-Null
